# Refactor synthetic array:
#  - insert a new "statut_name" column (C) between statut_label (B) and NCTId (C)
#  - derive its text from the existing statut_label (B) column
#  - the old "results_1y" column (I) is dropped/replaced by a new synthetic
#    placeholder column (now J) that is always FALSE
#  - results_3y / results / intervention_type simply shift one column to the
#    right (J->K, K->L, L->M), which Columns.Insert already does for us

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row on the sheet (header + data rows).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# 1) Insert a blank column at C - this shifts old C..L to D..M, and
#    carries the header style (s="1") along with it automatically.
$ws.Columns("C:C").Insert()

# 2) Header for the new column.
$ws.Cells.Item(1, 3).Value2 = "statut_name"

# 3) Map each row's statut_label (column B) to the new statut_name text
#    (column C), and reset the new results_1y column (now J) to FALSE.
$labelToName = @{
    "noir"   = "pas de résultat ni de publication"
    "rouge"  = "résultat et / ou publication posté"
    "orange" = "résultat et / ou publication posté dans les 36 mois"
}

for ($r = 2; $r -le $lastRow; $r++) {
    $label = $ws.Cells.Item($r, 2).Value2
    if ($labelToName.ContainsKey($label)) {
        $ws.Cells.Item($r, 3).Value2 = $labelToName[$label]
    }
    # Column J (results_1y) is the synthetic/dropped column - always FALSE now.
    $ws.Cells.Item($r, 10).Value2 = $false
}
